$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Table UI test cases appended below the existing testData rows
$ws.Range("A6").Value = "TEST-34227"
$ws.Range("A7").Value = "TEST-34229"

# Leave the selection on the next empty row, as the author's session did
$ws.Range("A8").Select()
